# Update the metrics table: refresh the rstd (and related) values for the
# existing models, drop the 3 worst models (model_2_0_22/23/24), and shrink
# the sheet dimension accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, "model_2_0_0", -0.1116691924315629, -1.168487262354351, 0.01095464815799552, -0.4044515863559262, 1.230290293693542, 1.338636517524719, 0.5054377317428589, 0.9465429782867432),
    @(3, "model_2_0_1", -0.03072692416193235, -0.8895490735474945, -0.02068650589514398, -0.280480688069662, 1.140710949897766, 1.166444301605225, 0.5216074585914612, 0.8629916310310364),
    @(4, "model_2_0_2", 0.02789344608488642, -0.4850405088568757, -0.2989322279059716, -0.1836141138045986, 1.075835466384888, 0.9167355298995972, 0.6638010740280151, 0.7977074980735779),
    @(5, "model_2_0_3", 0.03490790739239602, -0.4996354503766627, -0.2714205842754098, -0.1808746918088424, 1.06807267665863, 0.9257452487945557, 0.6497415900230408, 0.7958613038063049),
    @(6, "model_2_0_4", 0.05048773079356439, -0.4781944518018022, -0.2472191346941528, -0.1618416210342182, 1.050830245018005, 0.9125095009803772, 0.6373738050460815, 0.7830337882041931),
    @(7, "model_2_0_5", 0.0795293074286908, -0.4673792184961558, -0.1790631184230553, -0.1322772721412313, 1.018689870834351, 0.9058330059051514, 0.6025435924530029, 0.7631086111068726),
    @(8, "model_2_0_6", 0.09128341483821201, -0.5743596362407899, -0.129238178085513, -0.1663746861103768, 1.005681395530701, 0.9718735218048096, 0.5770812630653381, 0.7860888242721558),
    @(9, "model_2_0_7", 0.1035518332409231, -0.54195770083781, -0.0991243681916818, -0.1399171106663188, 0.9921039938926697, 0.9518713355064392, 0.5616920590400696, 0.7682575583457947),
    @(10, "model_2_0_8", 0.1048417575680565, -0.5475431704420253, -0.09111689224449848, -0.1397682888418885, 0.9906764030456543, 0.9553192853927612, 0.5575999021530151, 0.7681572437286377),
    @(11, "model_2_0_9", 0.1051498466166808, -0.548298422195068, -0.0889979756492536, -0.1393785596194077, 0.9903354644775391, 0.9557855129241943, 0.5565170645713806, 0.7678945660591125),
    @(12, "model_2_0_10", 0.1055527487061253, -0.544213914595415, -0.08935958452450166, -0.1375268953642821, 0.9898895621299744, 0.9532641172409058, 0.5567018389701843, 0.7666466236114502),
    @(13, "model_2_0_12", 0.1057828487364876, -0.5431030860147146, -0.08762109535155238, -0.1363679234409256, 0.9896349310874939, 0.9525784254074097, 0.5558134317398071, 0.7658655643463135),
    @(14, "model_2_0_11", 0.105877135790324, -0.5417579689820391, -0.08801320061981244, -0.135855547007101, 0.9895305633544922, 0.9517481327056885, 0.556013822555542, 0.7655202150344849),
    @(15, "model_2_0_13", 0.1061101087350805, -0.5388604121024672, -0.08848636580733271, -0.134619223556312, 0.9892727136611938, 0.9499593377113342, 0.556255578994751, 0.7646870613098145),
    @(16, "model_2_0_14", 0.1061197239702412, -0.5388051029080809, -0.08842507915592179, -0.1345705566833553, 0.9892620444297791, 0.9499253034591675, 0.5562242865562439, 0.7646541595458984),
    @(17, "model_2_0_15", 0.1063548692600413, -0.5369043855157176, -0.08749645086726399, -0.1333174625567273, 0.9890019297599792, 0.9487518668174744, 0.5557497143745422, 0.7638096213340759),
    @(18, "model_2_0_16", 0.1064499290018753, -0.5363708606904387, -0.0868037926111076, -0.1328116483704282, 0.9888966679573059, 0.9484225511550903, 0.555395781993866, 0.7634687423706055),
    @(19, "model_2_0_19", 0.1065233054915827, -0.5355036481673132, -0.08688596921883529, -0.1324204620900151, 0.9888154268264771, 0.947887122631073, 0.5554378032684326, 0.7632050514221191),
    @(20, "model_2_0_20", 0.1065233054915827, -0.5355036481673132, -0.08688596921883529, -0.1324204620900151, 0.9888154268264771, 0.947887122631073, 0.5554378032684326, 0.7632050514221191),
    @(21, "model_2_0_21", 0.1065233054915827, -0.5355036481673132, -0.08688596921883529, -0.1324204620900151, 0.9888154268264771, 0.947887122631073, 0.5554378032684326, 0.7632050514221191),
    @(22, "model_2_0_18", 0.1065234553996213, -0.5355025330923284, -0.08688407504000861, -0.1324192123935768, 0.9888152480125427, 0.94788658618927, 0.5554367899894714, 0.7632042765617371),
    @(23, "model_2_0_17", 0.1065385361164388, -0.535304967687114, -0.08692146679067658, -0.1323368537891831, 0.9887985587120056, 0.94776451587677, 0.5554559230804443, 0.7631487846374512)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
    $ws.Cells.Item($r, 9).Value = $row[9]
}

# The old table had 3 extra, now-dropped models at the bottom (rows 24-26).
$ws.Rows("24:26").Delete()
